# Update BAU Dispatch Priority by Electricity Source so that all variable
# renewable sources (onshore wind, solar PV, solar thermal, geothermal,
# offshore wind) are assigned Priority 1, while natural gas peaker and
# petroleum move to Priority 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDPbES")

# Row 8  -> solar thermal   : 2 -> 1
# Row 10 -> geothermal      : 2 -> 1
# Row 11 -> petroleum       : 1 -> 2
# Row 12 -> natural gas peaker : 1 -> 2
# Row 14 -> offshore wind   : 2 -> 1
$ws.Range("B8").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 2
$ws.Range("B14").Value = 1

# Match the author's final cursor/selection position on the sheet.
$ws.Range("B13").Select()

$wb.Save()
